$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("elite swimmers times during Olympic years", $true, $false, $false, $false, $false, $true, 1, $false, "elite swimmers" + [char]0x2019 + " times during Olympic years", 2)
